$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.328.24"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.640.43"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.40"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.91"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.648.35"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.84%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.105.18"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.305.17"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.02"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.650.01"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "349.01"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.54%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.23"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.25"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.415"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.14%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.14"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.45"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.90"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.58"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.64"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.941"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -12.14%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.861"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.16%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.66"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0992"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "277.43"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.600"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.56"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.067.73"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.48%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.75%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.31"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.71"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.43%  "
